$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# Insert a brand-new slide ("Lambda Function" / "Useful for short, ...")
# at position 12 (1-based), pushing the existing slide 12 ("Lambda
# Function" / "typical form of a lambda function") down to 13 and the
# old slide 13 ("Lambda Function" / "can capture local variables") to 14.
# ------------------------------------------------------------------

$refLayout = $p.Slides.Item(12).Layout
$news = $p.Slides.Add(12, $refLayout)

# ---- Title placeholder -------------------------------------------------
$title = $news.Shapes.Item(1)
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Lambda Function"
$titleTr.Font.Bold = $true
$titleTr.Font.Underline = $true

# ---- Body / content placeholder ----------------------------------------
$body = $news.Shapes.Item(2)

# Reposition/resize to match the target layout override.
# (Shape.Left/Top/Width/Height are expressed in points; 1 pt = 12700 EMU.)
$body.Left = 395536 / 12700
$body.Top = 1600200 / 12700
$body.Width = 8507288 / 12700
$body.Height = 4525963 / 12700

$tf = $body.TextFrame
$tr = $tf.TextRange

$line1 = "Useful for short, temporary, or one-time-use functions, especially when working with"
$line2 = "STL algorithms (e.g. sort, transform, etc.)"
$line3 = "Threads"
$line4 = "Callbacks or event handlers"

# Build paragraph by paragraph via InsertAfter (rather than one Text
# assignment containing embedded carriage returns) so every paragraph
# keeps a proper lang="en-US" run property.
$tr.Text = $line1
[void]$body.TextFrame.TextRange.InsertAfter("`r" + $line2)
[void]$body.TextFrame.TextRange.InsertAfter("`r" + $line3)
[void]$body.TextFrame.TextRange.InsertAfter("`r" + $line4)

$full = $body.TextFrame.TextRange

# Sub-bullets (paragraphs 2-4) are one indent level deeper.
$full.Paragraphs(2, 1).IndentLevel = 2
$full.Paragraphs(3, 1).IndentLevel = 2
$full.Paragraphs(4, 1).IndentLevel = 2

# Bold + red highlights inside the first line.
$redBold = @(
    @(12, 5),    # "short"
    @(19, 9),    # "temporary"
    @(33, 22)    # "one-time-use functions"
)
foreach ($pair in $redBold) {
    $run = $full.Characters($pair[0], $pair[1])
    $run.Font.Bold = $true
    $run.Font.Color.RGB = 255
}

Write-Host "Inserted slide at index 12; total slides:" $p.Slides.Count
